$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row above row 41 ("wohnhaft Sachsen?"); this pushes rows
# 41..48 down to 42..49 and carries their styles/row-heights/validations
# along automatically.
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with the new question text, mirroring the look
# of row 40 ("Minijob?") which uses the plain integer-formatted style.
$ws.Range("A41").Value = "juenger als 23 oder vor 1940 geboren?"
$ws.Range("B40").Copy() | Out-Null
$ws.Range("B41").PasteSpecial(-4122) | Out-Null
$ws.Range("B41").ClearContents()
$excel.CutCopyMode = 0

# Match the author's final selection/view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("A39").Select()
